# heliaphen_experiments: update the "15HP04" experiment's start/end dates
# (row 6, columns C/D) and leave the selection on the last-edited cell,
# matching the author's manual edit in Excel/Calc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "2015-06-03"
$ws.Range("D6").Value = "2015-06-12"

$ws.Range("D6").Select()
